# Add the missing "Log out on" time-stamp cells in column X of Sheet1.
# Rows 5, 7, 10 and 12 currently have an empty X cell (logout not recorded
# yet); the author typed the literal text "NULL" into each of them.
#
# X5 already carried the column's usual "time" style (s="5" -> numFmtId 20,
# centered, thin border) even while empty, but X7/X10/X12 were still using
# the plain bordered/centered style (s="2") left over from when the column
# was first formatted. Typing "NULL" into them also picked up the format
# already used elsewhere in the same column (as seen on X6/X8/X9/X11), so we
# copy that cell's formatting across after setting the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X5").Value = "NULL"
$ws.Range("X7").Value = "NULL"
$ws.Range("X10").Value = "NULL"
$ws.Range("X12").Value = "NULL"

# Bring X7/X10/X12 in line with the time-format style already used by the
# rest of column X (e.g. X6), matching the author's formatting.
$ws.Range("X6").Copy()
$ws.Range("X7").PasteSpecial(-4122)
$ws.Range("X6").Copy()
$ws.Range("X10").PasteSpecial(-4122)
$ws.Range("X6").Copy()
$ws.Range("X12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Final selection left on X12, matching the saved workbook state.
$ws.Range("X12").Select()
